# ===========================================================================
# Applies the OOPSpresentation.pptx edit:
#  1) Delete the empty picture-only slide (old slide 8 - "Content Placeholder
#     3" with no text, layout "Blank").
#  2) Insert a brand-new "53 years ago..." slide right after slide 1.
#  3) Append a brand-new "class" (Modifiers/Class name/...) slide at the end.
#  4) Small in-place text tweaks on the "Objects" and "A class is a
#     blueprint..." slides (run clean-up with identical visible text).
# ===========================================================================

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Step 1: remove the old picture-only slide (Content Placeholder 3, Blank
# layout, no text) -- it sits right after "Differences..."/"PROCEDURAL
# LANGUAGES"/"OBJECT ORIENTED LANGUAGES"/"Characteristics of oops" slides.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $cand = $p.Slides.Item($i)
    if ($cand.Shapes.Count -eq 1) {
        $cs = $cand.Shapes.Item(1)
        if ($cs.Name -eq "Content Placeholder 3" -and $cs.HasTextFrame -eq $false) {
            $cand.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------------
# Step 2: insert the new "53 years ago..." slide as slide 2 (right after the
# title slide "OOPS").
# ---------------------------------------------------------------------------
$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)
$s2 = $p.Slides.AddSlide(2, $titleAndContent)

$s2.Shapes.Item(1).TextFrame.TextRange.Text = "53 years ago…"

$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "Alan Kay in 1967 coined " + [char]8220 + "object oriented programming" + [char]8221 + ".`rInspired from cells .`rDeveloped Smalltalk at Xerox.`rAccording to Alan Kay, the essential ingredients of OOP are:`rMessage passing`rEncapsulation`rDynamic binding`r`r"

# "According" -> split "A" / "ccording " as two runs (matches target XML).
$full2 = $body2.Text
$idx2 = $full2.IndexOf("According to Alan Kay")
$null = $body2.Characters($idx2 + 1, 1).Text
$body2.Characters($idx2 + 2, 9).Text = "ccording "

# ---------------------------------------------------------------------------
# Step 3: append the new "class" slide (Modifiers / Class name / Superclass /
# Interfaces / Body) at the very end, mirroring slide13's "Blank" layout +
# free-floating rectangle textbox style.
# ---------------------------------------------------------------------------
$titleOnly = $p.SlideMaster.CustomLayouts.Item(6)
$sNew = $p.Slides.AddSlide($p.Slides.Count + 1, $titleOnly)

$titlePh = $sNew.Shapes.Placeholders.Item(1)
$titlePh.TextFrame.TextRange.Text = "class"

$rect = $sNew.Shapes.AddTextbox(1, 1455313, 1700011, 9465972, 3517438)
$rect.Name = "Rectangle 2"
$rtf = $rect.TextFrame
$rtf.WordWrap = $true
$rtf.AutoSize = 1

$rtr = $rtf.TextRange
$rtr.Text = "Modifiers: A class can be public or has default access `r`r`rClass name: The name should begin with a initial letter `r`rSuperclass(if any): The name of the class" + [char]8217 + "s parent (superclass), if any, preceded by the keyword extends. A class can only extend (subclass) one parent.`r`rInterfaces(if any): A comma-separated list of interfaces implemented by the class, if any, preceded by the keyword implements. A class can implement more than one interface.`r`r`rBody: The class body surrounded by braces, { }."

$fullr = $rtr.Text
$rtr.Font.Size = 16

# Bold the labels at the start of each relevant paragraph.
$pos = $fullr.IndexOf("Modifiers:")
$rtr.Characters($pos + 1, 9).Font.Bold = $true

$pos = $fullr.IndexOf("Class name:")
$rtr.Characters($pos + 1, 11).Font.Bold = $true

$pos = $fullr.IndexOf("Superclass(if any):")
$rtr.Characters($pos + 1, 20).Font.Bold = $true

$pos = $fullr.IndexOf("Interfaces(if any):")
$rtr.Characters($pos + 1, 20).Font.Bold = $true

$pos = $fullr.IndexOf("Body:")
$rtr.Characters($pos + 1, 5).Font.Bold = $true

# ---------------------------------------------------------------------------
# Step 4: minor, visually-invisible text clean-up (run merges) on the
# "Objects" slide and the "A class is a blueprint..." slide.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $cand = $p.Slides.Item($i)
    if ($cand.Shapes.Count -eq 2 -and $cand.Shapes.Item(1).HasTextFrame) {
        if ($cand.Shapes.Item(1).TextFrame.TextRange.Text -eq "Objects") {
            $tr = $cand.Shapes.Item(2).TextFrame.TextRange
            $full = $tr.Text
            $fidx = $full.IndexOf(": It is represented by attributes of an object.")
            if ($fidx -ge 0) {
                $tr.Characters($fidx + 1, 49).Text = ": It is represented by attributes of an object.  "
            }
            break
        }
    }
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $cand = $p.Slides.Item($i)
    if ($cand.Shapes.Count -eq 1 -and $cand.Shapes.Item(1).HasTextFrame) {
        $t0 = $cand.Shapes.Item(1).TextFrame.TextRange.Text
        if ($t0.StartsWith("A") -and $t0.Contains("blueprint")) {
            $tr = $cand.Shapes.Item(1).TextFrame.TextRange
            $tr.Characters(1, 2).Text = "A "
            break
        }
    }
}
